$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New amount-spent figures for Nandan (rows already existed, just add G values) ---
$ws.Range("G2").Value = 1200
$ws.Range("G3").Value = 1200

# --- New note next to Nandan's amount spent ---
$ws.Range("H4").Value = "5204-2400=2804"

# --- New GST breakup table starting at row 22 ---
$ws.Range("A22").Value = "18th Aug 2021"
$ws.Range("D22").Value = "No gst"

$ws.Range("A23").Value = "Adapter 5v 2.4 Amp"
$ws.Range("B23").Value = 1
$ws.Range("C23").Value = 296.61
$ws.Range("D23").Value = 296.61

$ws.Range("A24").Value = "Memory Card 32gb"
$ws.Range("B24").Value = 1
$ws.Range("C24").Value = 466.11
$ws.Range("D24").Value = 466.11

$ws.Range("A25").Value = "BreadBoard"
$ws.Range("B25").Value = 1
$ws.Range("C25").Value = 67.8
$ws.Range("D25").Value = 67.8

$ws.Range("A26").Value = "RaspberryPi 3b"
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = 2796.61
$ws.Range("D26").Value = 2796.61

$ws.Range("A27").Value = "Pi Camera"
$ws.Range("B27").Value = 1
$ws.Range("C27").Value = 466.11
$ws.Range("D27").Value = 466.11

$ws.Range("A28").Value = "SG90 Servo motor"
$ws.Range("B28").Value = 1
$ws.Range("C28").Value = 152.55
$ws.Range("D28").Value = 152.55

$ws.Range("A29").Value = "Mini Buzzer"
$ws.Range("B29").Value = 1
$ws.Range("C29").Value = 16.95
$ws.Range("D29").Value = 16.95

$ws.Range("A30").Value = "LED"
$ws.Range("B30").Value = 2
$ws.Range("C30").Value = 1.69
$ws.Range("D30").Value = 3.38

$ws.Range("A31").Value = "Ultrasonic Sensor HC-SR04"
$ws.Range("B31").Value = 1
$ws.Range("C31").Value = 101.7
$ws.Range("D31").Value = 101.7

$ws.Range("A32").Value = "1/4 Resistor"
$ws.Range("B32").Value = 4
$ws.Range("C32").Value = 0.85
$ws.Range("D32").Value = 3.4

$ws.Range("A33").Value = "Jumper Wire"
$ws.Range("B33").Value = 15
$ws.Range("C33").Value = 2.55
$ws.Range("D33").Value = 38.25

$ws.Range("A34").Value = "Total Without GST"
$ws.Range("D34").Formula = "=SUM(D23:D33)"

$ws.Range("C35").Value = "CGST"
$ws.Range("D35").Value = 396.85

$ws.Range("C36").Value = "SGST"
$ws.Range("D36").Value = 396.85

$ws.Range("A37").Value = "Total"
$ws.Range("D37").Formula = "=SUM(D34:D36)"

# --- Styling: match header/total formatting used elsewhere in the sheet ---
$ws.Range("A22").Font.Bold = $true
$ws.Range("A22").Font.Underline = $true
$ws.Range("D22").Font.Bold = $true
$ws.Range("A37").Font.Bold = $true
$ws.Range("D37").Font.Bold = $true

# --- View state: scroll so column B is leftmost and selection sits on H6 ---
$ws.Range("H6").Select()
$excel.ActiveWindow.ScrollColumn = 2
